# Update countries & provincias Spain
#
# "España" moves up in the country ranking (now ranks ahead of Colombia and
# Argentina), and the country stats table + "last updated" timestamp are
# refreshed with newer data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Last updated timestamp (row 1) ---------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Octubre de 2020 a las 20:42"

# --- Country ranking reshuffle: España jumps ahead of Colombia/Argentina --
# Row 8 was Colombia, row 9 was Argentina, row 10 was España.
# After the update: row 8 = España (new data), row 9 = Colombia (old row 8
# data), row 10 = Argentina (old row 9 data).
$ws.Cells.Item(8, 1).Value = "España"
$ws.Cells.Item(9, 1).Value = "Colombia"
$ws.Cells.Item(10, 1).Value = "Argentina"

# Column numbers: B=2 Casos totales, C=3 Nuevos casos, D=4 Casos activos,
# E=5 Recuperados, F=6 Casos criticos, G=7 Muertes hoy, H=8 Muertes

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 8010811
$ws.Cells.Item(4, 3).Value = 18813
$ws.Cells.Item(4, 4).Value = 5141007
$ws.Cells.Item(4, 5).Value = 2649960
$ws.Cells.Item(4, 7).Value = 149
$ws.Cells.Item(4, 8).Value = 219844

# Row 8: España (new, updated data)
$ws.Cells.Item(8, 2).Value = 918223
$ws.Cells.Item(8, 3).Value = 9286
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 7).Value = 65
$ws.Cells.Item(8, 8).Value = 33124

# Row 9: Colombia
$ws.Cells.Item(9, 2).Value = 911316
$ws.Cells.Item(9, 4).Value = 789787
$ws.Cells.Item(9, 5).Value = 93695
$ws.Cells.Item(9, 8).Value = 27834

# Row 10: Argentina
$ws.Cells.Item(10, 2).Value = 894206
$ws.Cells.Item(10, 4).Value = 721380
$ws.Cells.Item(10, 5).Value = 148958
$ws.Cells.Item(10, 8).Value = 23868

# Row 13: Peru
$ws.Cells.Item(13, 2).Value = 743479
$ws.Cells.Item(13, 3).Value = 8505
$ws.Cells.Item(13, 5).Value = 609872
$ws.Cells.Item(13, 7).Value = 96
$ws.Cells.Item(13, 8).Value = 32779

# Row 25
$ws.Cells.Item(25, 2).Value = 329510
$ws.Cells.Item(25, 3).Value = 3219
$ws.Cells.Item(25, 5).Value = 42899

# Row 29
$ws.Cells.Item(29, 2).Value = 182708
$ws.Cells.Item(29, 3).Value = 844
$ws.Cells.Item(29, 4).Value = 154238

# Row 53
$ws.Cells.Item(53, 2).Value = 85136
$ws.Cells.Item(53, 3).Value = 841
$ws.Cells.Item(53, 4).Value = 38904
$ws.Cells.Item(53, 5).Value = 44931
$ws.Cells.Item(53, 7).Value = 14
$ws.Cells.Item(53, 8).Value = 1301

# Row 72
$ws.Cells.Item(72, 2).Value = 43351
$ws.Cells.Item(72, 3).Value = 823
$ws.Cells.Item(72, 5).Value = 18160
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 1827

# Row 109
$ws.Cells.Item(109, 2).Value = 10088
$ws.Cells.Item(109, 3).Value = 87
$ws.Cells.Item(109, 4).Value = 7550
$ws.Cells.Item(109, 5).Value = 2466
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 72

# Row 124
$ws.Cells.Item(124, 2).Value = 5683
$ws.Cells.Item(124, 3).Value = 14
$ws.Cells.Item(124, 4).Value = 5322
$ws.Cells.Item(124, 5).Value = 248

# Row 138
$ws.Cells.Item(138, 2).Value = 4197
$ws.Cells.Item(138, 3).Value = 3
$ws.Cells.Item(138, 4).Value = 3790
$ws.Cells.Item(138, 5).Value = 374
